$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 37, shifting existing rows 37-40 down to 38-41.
$ws.Rows("37").Insert()

# Populate the newly inserted row 37 with the new weekly record.
$ws.Range("A37").Value = 11
$ws.Range("B37").Value = "Vega Monumental Concepción"
$ws.Range("C37").Value = "Bíobío"
$ws.Range("D37").Value = 44722
$ws.Range("E37").Value = 8
$ws.Range("F37").Value = 100112013
$ws.Range("G37").Value = "Alcachofa"
$ws.Range("H37").Value = "Española"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 120
$ws.Range("K37").Value = 20000
$ws.Range("L37").Value = 22000
$ws.Range("M37").Value = 21333
$ws.Range("N37").Value = "$/caja 30 unidades"
$ws.Range("O37").Value = "Provincia de Limarí"
$ws.Range("P37").Value = 711
$ws.Range("Q37").Value = 30
$ws.Range("R37").Value = "Hortaliza"
